$d = $word.ActiveDocument

# Locate the last answer of the "pick hand tap" question:
#   "Consiste nel martellare le corde sulla tastiera con una mano qualsiasi"
$rng = $d.Content
$rng.Find.Execute("Consiste nel martellare le corde sulla tastiera con una mano qualsiasi", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$splitPos = $rng.Start

# Type the new answer text right after the existing run. Word glues the
# hidden "_GoBack" bookmark to the end of whatever was just typed, so it
# will correctly end up after this new text once the paragraph is split
# below.
$rng.Text = " Composizione accordo random"

# Insert a paragraph break *before* the text we just typed, turning
# "Composizione accordo random" into its own new list paragraph.
$breakRng = $d.Range($splitPos, $splitPos)
$breakRng.InsertParagraphAfter()

# Find the freshly created paragraph that holds the new text and promote it
# from ilvl 1 (inherited from the answer it was split off of) to ilvl 0,
# i.e. a top-level question bullet, matching its sibling questions.
$newParaIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq " Composizione accordo random") {
        $newParaIdx = $i
        break
    }
}
$newPara = $d.Paragraphs.Item($newParaIdx)
$newPara.Range.ListFormat.ListLevelNumber = 1

# Insert a new blank paragraph right after it, styled like the blank line
# that separates every instrument's quiz section (bold, 15pt).
$newPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Item($newParaIdx + 1)
$blankPara.Style = "Normal"
$blankPara.Range.Font.Bold = 1
$blankPara.Range.Font.Size = 15
$blankPara.Range.Font.SizeBi = 15

# InsertParagraphAfter() leaves a phantom empty run behind in the new blank
# paragraph; trim it so it matches its sibling blank paragraphs exactly
# (pPr only, no run).
$trimRng = $blankPara.Range.Duplicate
$trimRng.MoveEnd(1, -1) | Out-Null
$trimRng.Delete()
